$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value = "JD_018"
$ws.Cells.Item($row, 2).Value = "Senior UI Engineer"
$ws.Cells.Item($row, 3).Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 4
